$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the text of the existing activity note in F11 ---
$ws.Range("F11").Value = "Edit little about UI, wait partners answer, based on partner's design and work style suggested new plan and still discussion is going on"

# --- Fix a typo in the date value for row 10 ---
$ws.Range("A10").Value = 30.11

# --- Row 11: add the "Sum. Time" entry ---
$ws.Range("E11").Value = "1hr"

# --- Row 12 ---
$ws.Range("A12").Value = 2.12
$ws.Range("B12").Value = 0.54166666666666663
$ws.Range("B12").NumberFormat = "h:mm"
$ws.Range("C12").Value = 0.625
$ws.Range("C12").NumberFormat = "h:mm"
$ws.Range("E12").Value = "2hr"
$ws.Range("F12").Value = "Studied math conception more about it permutation with replacement and without replacement"

# --- Row 13 ---
$ws.Range("A13").Value = 3.12
$ws.Range("B13").Value = 0.95833333333333337
$ws.Range("B13").NumberFormat = "h:mm"
$ws.Range("C13").Value = 0.041666666666666664
$ws.Range("C13").NumberFormat = "h:mm"
$ws.Range("E13").Value = "2hr"
$ws.Range("F13").Value = "Check about function javascript"

# --- Row 14 ---
$ws.Range("A14").Value = 4.12
$ws.Range("B14").Value = 0.375
$ws.Range("B14").NumberFormat = "h:mm"
$ws.Range("C14").Value = 0.66666666666666663
$ws.Range("C14").NumberFormat = "h:mm"
$ws.Range("E14").Value = "7hr"
$ws.Range("F14").Value = "figure out second function(number system chart) and last function first try out and still figure out first one"

# --- Row 15 ---
$ws.Range("A15").Value = 5.12
$ws.Range("B15").Value = 0.375
$ws.Range("B15").NumberFormat = "h:mm"
$ws.Range("C15").Value = 0.625
$ws.Range("C15").NumberFormat = "h:mm"
$ws.Range("E15").Value = "6hr"
$ws.Range("F15").Value = "make other functions basics and try to change first function(number conversion)"

# --- Row 16 ---
$ws.Range("A16").Value = 6.12
$ws.Range("B16").Value = 0.70833333333333337
$ws.Range("B16").NumberFormat = "h:mm"
$ws.Range("C16").Value = 0.75
$ws.Range("C16").NumberFormat = "h:mm"
$ws.Range("E16").Value = "1hr"
$ws.Range("F16").Value = "organize function one html file"

# --- Update the active selected cell in the sheet view ---
$ws.Range("F17").Select()
